$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "31÷9=" "60÷7="
Replace-Text "44÷5=" "40÷5="
Replace-Text "40÷3=" "68÷3="
Replace-Text "35÷4=" "74÷2="
Replace-Text "13÷8=" "25÷7="
Replace-Text "18÷4=" "80÷8="
Replace-Text "47÷9=" "38÷9="
Replace-Text "15÷8=" "77÷3="
Replace-Text "17÷9=" "71÷8="
Replace-Text "72÷4=" "77÷3="
Replace-Text "31÷7=" "96÷8="
Replace-Text "99÷9=" "52÷4="
Replace-Text "94÷9=" "27÷7="
Replace-Text "80÷6=" "77÷3="
Replace-Text "64÷8=" "54÷2="
Replace-Text "65÷7=" "26÷3="
Replace-Text "19÷6=" "85÷7="
Replace-Text "47÷3=" "17÷2="
Replace-Text "48÷4=" "23÷9="
Replace-Text "25÷9=" "25÷6="
Replace-Text "20÷8=" "84÷8="
Replace-Text "15÷2=" "39÷9="
Replace-Text "44÷4=" "19÷2="
Replace-Text "12÷2=" "60÷5="
Replace-Text "33÷5=" "46÷6="

Write-Output "Done applying replacements"
